$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C12" = -10.984
    "C27" = -12.775
    "C32" = -13.676
    "C36" = -12.732
    "C38" = -12.703
    "C46" = -14.033
    "C54" = -13.137
    "C55" = -13.519
    "C56" = -13.314
    "C67" = -11.595
    "C69" = -11.038
    "C72" = -11.555
    "C83" = -13.283
    "C86" = -13.665
    "C91" = -10.98
    "C93" = -11.979
    "C99" = -12.715
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
